$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the commit diff.
# B/C columns are plain text (coin name / link) - direct assignment is safe.
# D column values are forced to Text format first so Excel does not
# auto-convert numeric-looking strings (e.g. '1.00' -> 1, '8.60' -> 8.6...),
# matching the inlineStr text cells in the original workbook.
# E column values (percentages with surrounding spaces) are never valid
# numbers, so they remain text automatically.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.357.25'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.220.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '110.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '289.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.83%  '
$ws.Range("E7").Value = '  -1.01%  '
$ws.Range("E8").Value = '  -0.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.596'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0909'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.29'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.60'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +13.19%  '
$ws.Range("E15").Value = '  -2.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.556.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.225.70'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.381.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.97%  '
$ws.Range("E21").Value = '  -2.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +13.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.40'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '234.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.78%  '
$ws.Range("E29").Value = '  -1.91%  '
$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '173.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '37.32'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -10.02%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.11'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '21.28'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0873'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.59'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.96'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.43%  '
$ws.Range("E37").Value = '  -1.86%  '
$ws.Range("E38").Value = '  -4.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0374'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.43%  '
$ws.Range("E40").Value = '  -1.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.31'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.229'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.10%  '
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.30'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.81%  '
$ws.Range("E46").Value = '  -2.44%  '
$ws.Range("E47").Value = '  -5.43%  '
$ws.Range("E48").Value = '  +1.50%  '
$ws.Range("E49").Value = '  +4.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '100.79'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.36%  '
